$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text in B2: remove "/RME" from the "26% S/LFM+CDM/RME/H:1" line
$ws.Range("B2").Value = "20% CR/LFINF+CDM/H:2`n9% CR/LFINF+CDL/HBET:3-5`n10% S+SL/LFM+CDL/H:1`n26% S/LFM+CDM/H:1`n15% CR+PC/LFM+CDL/H:1`n5% MUR/LWAL+CDN/H:1`n15% MCF/LWAL+CDL/H:1"

# Apply wrap text alignment to B2
$ws.Range("B2").WrapText = $true

# Set row height for row 2
$ws.Rows.Item(2).RowHeight = 395

# Set selection: B2:B11 (active cell lands on B2, the anchor corner)
$ws.Range("B2:B11").Select()
